# Insert two new weekly price rows for "Brócoli" (Segunda / Tercera calidad,
# fecha 2023-12-07) just above the existing row 666, pushing the rest of the
# table (old rows 666:697) down to 668:699.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A666:A667").EntireRow.Insert()

# New row 666: Brócoli, Segunda, fecha 45267 (2023-12-07)
$ws.Cells.Item(666, 1).Value = 1
$ws.Cells.Item(666, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(666, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(666, 4).Value = 45267
$ws.Cells.Item(666, 5).Value = 15
$ws.Cells.Item(666, 6).Value = 100112023
$ws.Cells.Item(666, 7).Value = "Brócoli"
$ws.Cells.Item(666, 8).Value = "Sin especificar"
$ws.Cells.Item(666, 9).Value = "Segunda"
$ws.Cells.Item(666, 10).Value = 1200
$ws.Cells.Item(666, 11).Value = 400
$ws.Cells.Item(666, 12).Value = 500
$ws.Cells.Item(666, 13).Value = 450
$ws.Cells.Item(666, 14).Value = '$/unidad'
$ws.Cells.Item(666, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(666, 16).Value = 450
$ws.Cells.Item(666, 17).Value = 1
$ws.Cells.Item(666, 18).Value = "Hortaliza"

# New row 667: Brócoli, Tercera, fecha 45267 (2023-12-07)
$ws.Cells.Item(667, 1).Value = 1
$ws.Cells.Item(667, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(667, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(667, 4).Value = 45267
$ws.Cells.Item(667, 5).Value = 15
$ws.Cells.Item(667, 6).Value = 100112023
$ws.Cells.Item(667, 7).Value = "Brócoli"
$ws.Cells.Item(667, 8).Value = "Sin especificar"
$ws.Cells.Item(667, 9).Value = "Tercera"
$ws.Cells.Item(667, 10).Value = 1200
$ws.Cells.Item(667, 11).Value = 300
$ws.Cells.Item(667, 12).Value = 350
$ws.Cells.Item(667, 13).Value = 325
$ws.Cells.Item(667, 14).Value = '$/unidad'
$ws.Cells.Item(667, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(667, 16).Value = 325
$ws.Cells.Item(667, 17).Value = 1
$ws.Cells.Item(667, 18).Value = "Hortaliza"
